$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1 matching the style of the existing header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the new column data for I and J (rows 2-11)
$iValues = @(1, 1, 1, 7, 1, 1, 1, 8, 4, 3)
$jValues = @(5, 6, 7, 8, 5, 6, 5, 8, 5, 4)

for ($r = 0; $r -lt 10; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
